$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last fully populated data row (row 35) down
# onto the two new rows (36 and 37) so the new entries render the same
# way (number formats, alignment, etc.) as the rest of the table.
$ws.Range("A35:E35").Copy()
$ws.Range("A36:E37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 36 : S.No 34 / Order 26377889 ---
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 26377889
$ws.Cells.Item(36, 3).Value = 42236
$ws.Cells.Item(36, 4).Value = 40000.11
$ws.Cells.Item(36, 5).Value = 44256

# --- Row 37 : S.No 35 / Order 26396131 ---
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 26396131
$ws.Cells.Item(37, 3).Value = 40018
$ws.Cells.Item(37, 4).Value = 37899.54
$ws.Cells.Item(37, 5).Value = 44257

# The F:I columns on rows 36/37 already contained the shared formulas
# (they previously evaluated against blank B36/B37, yielding ""); force
# each to recompute now that the rows have real data by re-applying the
# existing formula text.
$ws.Cells.Item(36, 6).Formula = $ws.Cells.Item(36, 6).Formula
$ws.Cells.Item(36, 7).Formula = $ws.Cells.Item(36, 7).Formula
$ws.Cells.Item(36, 8).Formula = $ws.Cells.Item(36, 8).Formula
$ws.Cells.Item(36, 9).Formula = $ws.Cells.Item(36, 9).Formula

$ws.Cells.Item(37, 6).Formula = $ws.Cells.Item(37, 6).Formula
$ws.Cells.Item(37, 7).Formula = $ws.Cells.Item(37, 7).Formula
$ws.Cells.Item(37, 8).Formula = $ws.Cells.Item(37, 8).Formula
$ws.Cells.Item(37, 9).Formula = $ws.Cells.Item(37, 9).Formula

$excel.CalculateFull()

# Move the view: scroll the frozen pane back to the top of the data and
# select I3, matching where the author left off after entering the data.
$ws.Range("I3").Select()

$wb.Save()
